$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 232, shifting existing rows 232-271 down to 233-272
$ws.Rows("232:232").Insert()

# Populate the new row 232 with the same constant/template values used by the
# surrounding rows, plus the new data values.
$ws.Range("A232").Value = 8
$ws.Range("B232").Value = "Terminal La Palmera de La Serena"
$ws.Range("C232").Value = "Coquimbo"
$ws.Range("D232").Value = 45258
$ws.Range("D232").NumberFormat = $ws.Range("D233").NumberFormat
$ws.Range("E232").Value = 4
$ws.Range("F232").Value = 100112040
$ws.Range("G232").Value = "Cilantro"
$ws.Range("H232").Value = "Sin especificar"
$ws.Range("I232").Value = "Primera"
$ws.Range("J232").Value = 2400
$ws.Range("K232").Value = 2300
$ws.Range("L232").Value = 2500
$ws.Range("M232").Value = 2400
$ws.Range("N232").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O232").Value = "Provincia del Elquí"
$ws.Range("P232").Value = 1600
$ws.Range("Q232").Value = 1.5
$ws.Range("R232").Value = "Hortaliza"
